$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add new "Label" header column (H), matching the style of the other headers ---
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)  # xlPasteFormats
$ws.Cells.Item(1, 8).Value = "Label"

# --- Update slightly-changed prediction/error values (re-fit results) ---
$ws.Range("D3").Value = 0.5004278364299665
$ws.Range("E3").Value = 0.5004278364299665

$ws.Range("D4").Value = 0.3190943878203306
$ws.Range("E4").Value = 0.3190943878203306

$ws.Range("D9").Value = 0.6432560361952234
$ws.Range("E9").Value = 0.3567439638047766

$ws.Range("D11").Value = 0.1152678792958895
$ws.Range("E11").Value = 0.8847321207041106
$ws.Range("F11").Value = 125.4838027954102

# --- Populate new Label column (H2:H21) ---
# Rows 2-6 and 12-16 correspond to Control patients -> Label 0
# Rows 7-11 and 17-21 correspond to MDD patients -> Label 1
$labels = @(0, 0, 0, 0, 0, 1, 1, 1, 1, 1, 0, 0, 0, 0, 0, 1, 1, 1, 1, 1)
for ($i = 0; $i -lt $labels.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 8).Value = $labels[$i]
}
